# Restore cell C10 on the "Rules" sheet to value 1 (was 18),
# per revision #b91add373b8fd054b57e7f5ed6615a705027d31a.TEST restore/save.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1.0
